$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily records (date-serial, nuovi pos., somma mobile 7gg., somma mobile 7gg. per 100mila ab.)
# aggiornamento fino a 1/09/2021
$newData = @(
    @(44432, 1, 10, 117.827265229174),
    @(44433, 2, 12, 141.3927182750088),
    @(44434, 1, 11, 129.6099917520914),
    @(44435, 0, 11, 129.6099917520914),
    @(44436, 4, 15, 176.740897843761),
    @(44437, 5, 13, 153.1754447979262),
    @(44438, 3, 16, 188.5236243666784),
    @(44439, 0, 15, 176.740897843761),
    @(44440, 0, 13, 153.1754447979262)
)

$lastRow = 357
$ws.Range("A$lastRow`:D$lastRow").Copy()

foreach ($rec in $newData) {
    $lastRow = $lastRow + 1
    $ws.Range("A$lastRow`:D$lastRow").PasteSpecial(-4122)
    $ws.Cells.Item($lastRow, 1).Value = $rec[0]
    $ws.Cells.Item($lastRow, 2).Value = $rec[1]
    $ws.Cells.Item($lastRow, 3).Value = $rec[2]
    $ws.Cells.Item($lastRow, 4).Value = $rec[3]
}
